$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 8 (shifts existing rows 8+ down by one)
$ws.Rows.Item(8).Insert()

# Copy the formatting of the row below (the original row 8 content, now at
# row 9, which carries the "boolean attribute row" styling) onto the newly
# inserted blank row 8, so it matches the existing TRUE/FALSE attribute rows.
$ws.Range("A9:AG9").Copy()
$ws.Range("A8:AG8").PasteSpecial(-4122)  # xlPasteFormats

# Set the new attribute row's label and values
$ws.Range("A8").Value = "Force"
$ws.Range("B8:AG8").Value = $false

# Rebuild data validation ranges so they cover the newly inserted row.
$ws.Cells.Validation.Delete()
$ws.Range("A7:A9").Validation.Add(0, 1, 1, "")
$ws.Range("B7:AG9").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# Re-freeze the header panes one row lower (B11 instead of B10) to account
# for the newly inserted row, and leave the selection on A9 as in the
# saved file.
[void]($excel.ActiveWindow.FreezePanes = $false)
[void]($ws.Range("B11").Select())
[void]($excel.ActiveWindow.FreezePanes = $true)
[void]($ws.Range("A9").Select())

Write-Host "Done"
